$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4257.542554456429,
    4136.71208971015,
    4136.71208971015,
    3993.298809462478,
    3993.298809462478,
    3993.298809462478,
    3917.468736119041,
    3917.468736119041,
    3917.468736119041,
    3834.185813160205,
    3815.950759670592
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
